# The workbook tracks daily/weekly price observations for "Perejil" (parsley)
# at "Feria Lagunitas de Puerto Montt". A new weekly record is inserted as a
# new row 312 (pushing the former rows 312..373 down to 313..374).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 312; Excel shifts rows 312:373 down to
# 313:374 and keeps their values/styles intact.
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row 312 with the new observation.
$ws.Cells.Item(312, 1).Value  = 4
$ws.Cells.Item(312, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(312, 3).Value  = "Los Lagos"
$ws.Cells.Item(312, 4).Value  = 45015
$ws.Cells.Item(312, 5).Value  = 10
$ws.Cells.Item(312, 6).Value  = 100112044
$ws.Cells.Item(312, 7).Value  = "Perejil"
$ws.Cells.Item(312, 8).Value  = "Sin especificar"
$ws.Cells.Item(312, 9).Value  = "Primera"
$ws.Cells.Item(312, 10).Value = 50
$ws.Cells.Item(312, 11).Value = 6000
$ws.Cells.Item(312, 12).Value = 6000
$ws.Cells.Item(312, 13).Value = 6000
$ws.Cells.Item(312, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(312, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(312, 16).Value = 3000
$ws.Cells.Item(312, 17).Value = 2
$ws.Cells.Item(312, 18).Value = "Hortaliza"

# Column D (Fecha) should keep the same date/time number format the rest of
# the column uses.
$ws.Cells.Item(312, 4).NumberFormat = $ws.Cells.Item(313, 4).NumberFormat
